# Clear the "attached document" test data columns (H:L) for the data rows
# (2-5) of the Transmittals test sheet. The header row (row 1) and the
# rest of the row data are left untouched. Clearing these cells also
# leaves the related shared strings ("Document Register", "Test 1 ta.docx",
# "BrowseDocument.docx") unused so Excel drops them from sharedStrings.xml
# when it rewrites the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns H (AttachDocuments) through L (ReviewDocument), rows 2-5.
$ws.Range("H2:L5").ClearContents()

# Reset the view back to the top-left corner with a plain A1 selection,
# matching the simplified <sheetView> left by the edit.
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
